# Auto-generated script to update Leve profit-tracking sheets
# with refreshed market-price data (columns H-N) for specific rows,
# as produced by the scheduled market-data runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40: Stuck in the Moment
$ws.Range("H40").Value = 2082.5
$ws.Range("I40").Value = 1276.2
$ws.Range("J40").Value = 2274.476
$ws.Range("K40").Value = 1276.2
$ws.Range("L40").Value = 2274.476
$ws.Range("M40").Value = -1101.2
$ws.Range("N40").Value = -2624.476
# Row 44: Alive and Unwell
$ws.Range("H44").Value = 21179.2
$ws.Range("J44").Value = 21179.2
$ws.Range("L44").Value = 21179.2
$ws.Range("N44").Value = -22103.2
# Row 138: All-night Crafting
$ws.Range("H138").Value = 2175.2646
$ws.Range("I138").Value = 1205.2
$ws.Range("J138").Value = 3561.0715
$ws.Range("K138").Value = 3615.6
$ws.Range("L138").Value = 10683.2145
$ws.Range("M138").Value = 1524.4
$ws.Range("N138").Value = -20963.2145
# Row 141: Remedy for Reason
$ws.Range("H141").Value = 10060.179
$ws.Range("I141").Value = 1737.381
$ws.Range("J141").Value = 35028.57
$ws.Range("K141").Value = 5212.143
$ws.Range("L141").Value = 105085.71
$ws.Range("M141").Value = -32.14300000000003
$ws.Range("N141").Value = -115445.71

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 6512.0566
$ws.Range("I32").Value = 6204.061
$ws.Range("K32").Value = 6204.061
$ws.Range("M32").Value = -5917.061
# Row 63: Rivets Run through It
$ws.Range("H63").Value = 5516.4614
$ws.Range("I63").Value = 2713.4285
$ws.Range("J63").Value = 8786.666999999999
$ws.Range("K63").Value = 2713.4285
$ws.Range("L63").Value = 8786.666999999999
$ws.Range("M63").Value = -2027.4285
$ws.Range("N63").Value = -10158.667
# Row 66: A Riveting Revival (L)
$ws.Range("H66").Value = 5516.4614
$ws.Range("I66").Value = 2713.4285
$ws.Range("J66").Value = 8786.666999999999
$ws.Range("K66").Value = 13567.1425
$ws.Range("L66").Value = 43933.335
$ws.Range("M66").Value = -10135.1425
$ws.Range("N66").Value = -50797.335
# Row 74: As the Bolt Flies
$ws.Range("H74").Value = 1387.9375
$ws.Range("I74").Value = 1129.5555
$ws.Range("J74").Value = 2163.0833
$ws.Range("K74").Value = 1129.5555
$ws.Range("L74").Value = 2163.0833
$ws.Range("M74").Value = -255.5554999999999
$ws.Range("N74").Value = -3911.0833
# Row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value = 1387.9375
$ws.Range("I77").Value = 1129.5555
$ws.Range("J77").Value = 2163.0833
$ws.Range("K77").Value = 5647.7775
$ws.Range("L77").Value = 10815.4165
$ws.Range("M77").Value = -1279.7775
$ws.Range("N77").Value = -19551.4165

$ws = $wb.Worksheets.Item("BSM")
# Row 94: High Steal
$ws.Range("H94").Value = 932.4761999999999
$ws.Range("I94").Value = 733.4666999999999
$ws.Range("J94").Value = 1430
$ws.Range("K94").Value = 733.4666999999999
$ws.Range("L94").Value = 1430
$ws.Range("M94").Value = -282.4666999999999
$ws.Range("N94").Value = -2332

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found
$ws.Range("H31").Value = 7938779.5
$ws.Range("I31").Value = 1627.75
$ws.Range("J31").Value = 18521648
$ws.Range("K31").Value = 1627.75
$ws.Range("L31").Value = 18521648
$ws.Range("M31").Value = -1332.75
$ws.Range("N31").Value = -18522238
# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 7938779.5
$ws.Range("I34").Value = 1627.75
$ws.Range("J34").Value = 18521648
$ws.Range("K34").Value = 1627.75
$ws.Range("L34").Value = 18521648
$ws.Range("M34").Value = -1425.75
$ws.Range("N34").Value = -18522052
# Row 58: You Do the Heavy Lifting
$ws.Range("H58").Value = 1966.3793
$ws.Range("I58").Value = 1309.6666
$ws.Range("J58").Value = 3041
$ws.Range("K58").Value = 1309.6666
$ws.Range("L58").Value = 3041
$ws.Range("M58").Value = -1106.6666
$ws.Range("N58").Value = -3447
# Row 62: Splinter in the Sewers
$ws.Range("H62").Value = 27150.111
$ws.Range("I62").Value = 27087.75
$ws.Range("J62").Value = 27200
$ws.Range("K62").Value = 27087.75
$ws.Range("L62").Value = 27200
$ws.Range("M62").Value = -26463.75
$ws.Range("N62").Value = -28448
# Row 65: The Lumber of Their Discontent (L)
$ws.Range("H65").Value = 27150.111
$ws.Range("I65").Value = 27087.75
$ws.Range("J65").Value = 27200
$ws.Range("K65").Value = 135438.75
$ws.Range("L65").Value = 136000
$ws.Range("M65").Value = -132318.75
$ws.Range("N65").Value = -142240
# Row 132: Hull Lotta Damage
$ws.Range("H132").Value = 2545.348
$ws.Range("I132").Value = 1820.5714
$ws.Range("J132").Value = 3672.7778
$ws.Range("K132").Value = 5461.7142
$ws.Range("L132").Value = 11018.3334
$ws.Range("M132").Value = -2931.7142
$ws.Range("N132").Value = -16078.3334
# Row 134: Wood You Be Quiet
$ws.Range("H134").Value = 1596.25
$ws.Range("I134").Value = 899
$ws.Range("J134").Value = 1864.4231
$ws.Range("K134").Value = 2697
$ws.Range("L134").Value = 5593.2693
$ws.Range("M134").Value = -162
$ws.Range("N134").Value = -10663.2693
# Row 136: Turali Quality
$ws.Range("H136").Value = 1966.3793
$ws.Range("I136").Value = 1309.6666
$ws.Range("J136").Value = 3041
$ws.Range("K136").Value = 3928.9998
$ws.Range("L136").Value = 9123
$ws.Range("M136").Value = -1378.9998
$ws.Range("N136").Value = -14223

$ws = $wb.Worksheets.Item("CUL")
# Row 15: Pretty Enough to Eat
$ws.Range("H15").Value = 289.83334
$ws.Range("I15").Value = 247.8
$ws.Range("J15").Value = 500
$ws.Range("K15").Value = 743.4000000000001
$ws.Range("L15").Value = 1500
$ws.Range("M15").Value = -603.4000000000001
$ws.Range("N15").Value = -1780

$ws = $wb.Worksheets.Item("GSM")
# Row 57: Gold Is So Last Year
$ws.Range("H57").Value = 16500
$ws.Range("J57").Value = 23000
$ws.Range("L57").Value = 23000
$ws.Range("N57").Value = -24640
# Row 80: Needs More Prayerbell
$ws.Range("H80").Value = 2584.7058
$ws.Range("I80").Value = 2503.077
$ws.Range("J80").Value = 2850
$ws.Range("K80").Value = 2503.077
$ws.Range("L80").Value = 2850
$ws.Range("M80").Value = -1505.077
$ws.Range("N80").Value = -4846
# Row 83: With a Noise That Reaches Heaven (L)
$ws.Range("H83").Value = 2584.7058
$ws.Range("I83").Value = 2503.077
$ws.Range("J83").Value = 2850
$ws.Range("K83").Value = 12515.385
$ws.Range("L83").Value = 14250
$ws.Range("M83").Value = -7523.385000000002
$ws.Range("N83").Value = -24234
# Row 122: Awarding Academic Excellence
$ws.Range("H122").Value = 3392.8823
$ws.Range("J122").Value = 2972.25
$ws.Range("L122").Value = 8916.75
$ws.Range("N122").Value = -13816.75
# Row 132: On Board for Lar
$ws.Range("H132").Value = 696383.1
$ws.Range("I132").Value = 1303549.4
$ws.Range("J132").Value = 2478.8928
$ws.Range("K132").Value = 3910648.2
$ws.Range("L132").Value = 7436.678400000001
$ws.Range("M132").Value = -3908118.2
$ws.Range("N132").Value = -12496.6784

$ws = $wb.Worksheets.Item("LTW")
# Row 122: Hell on Leather
$ws.Range("H122").Value = 9900
$ws.Range("I122").Value = 10000
$ws.Range("J122").Value = 9800
$ws.Range("K122").Value = 30000
$ws.Range("L122").Value = 29400
$ws.Range("M122").Value = -27550
$ws.Range("N122").Value = -34300
# Row 132: Tenets of Tanning
$ws.Range("H132").Value = 27030068
$ws.Range("I132").Value = 35717012
$ws.Range("J132").Value = 4023.111
$ws.Range("K132").Value = 107151036
$ws.Range("L132").Value = 12069.333
$ws.Range("M132").Value = -107148506
$ws.Range("N132").Value = -17129.333

$ws = $wb.Worksheets.Item("WVR")
# Row 81: Where the Dragonflies, the Net Catches
$ws.Range("H81").Value = 2580
$ws.Range("I81").Value = 3575.5
$ws.Range("K81").Value = 7151
$ws.Range("M81").Value = -6090
# Row 84: To Kill a Dragon on Nameday (L)
$ws.Range("H84").Value = 2580
$ws.Range("I84").Value = 3575.5
$ws.Range("K84").Value = 35755
$ws.Range("M84").Value = -30451
# Row 107: Flax Wax
$ws.Range("H107").Value = 608.8889
$ws.Range("I107").Value = 613.3333
$ws.Range("K107").Value = 1839.9999
$ws.Range("M107").Value = 80.00009999999997
# Row 122: Heavy Armoire
$ws.Range("H122").Value = 3896.275
$ws.Range("I122").Value = 3077.6296
$ws.Range("J122").Value = 5596.5386
$ws.Range("K122").Value = 9232.888800000001
$ws.Range("L122").Value = 16789.6158
$ws.Range("M122").Value = -6782.888800000001
$ws.Range("N122").Value = -21689.6158
# Row 132: Comfy Cabins
$ws.Range("H132").Value = 2025.6066
$ws.Range("I132").Value = 1637.675
$ws.Range("J132").Value = 2764.524
$ws.Range("K132").Value = 4913.025
$ws.Range("L132").Value = 8293.572
$ws.Range("M132").Value = -2383.025
$ws.Range("N132").Value = -13353.572
# Row 136: Weaving the Envelope
$ws.Range("H136").Value = 5683171.5
$ws.Range("I136").Value = 7576409
$ws.Range("J136").Value = 3459.6365
$ws.Range("K136").Value = 22729227
$ws.Range("L136").Value = 10378.9095
$ws.Range("M136").Value = -22726677
$ws.Range("N136").Value = -15478.9095
